# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.389.51'
$ws.Range('E2').Value = '  +3.00%  '
$ws.Range('D3').Value = '2.627.35'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.06'
$ws.Range('E5').Value = '  +5.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.23'
$ws.Range('E6').Value = '  +1.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.607'
$ws.Range('E8').Value = '  +6.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.83'
$ws.Range('E9').Value = '  -2.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.103'
$ws.Range('E10').Value = '  +2.29%  '
$ws.Range('E11').Value = '  +6.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.342'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '3.091.81'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').Value = '60.309.60'
$ws.Range('E14').Value = '  +2.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.62'
$ws.Range('E15').Value = '  +3.37%  '
$ws.Range('D16').Value = '2.636.61'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('E17').Value = '  +1.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.58'
$ws.Range('E18').Value = '  +4.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '342.76'
$ws.Range('E19').Value = '  +2.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.36'
$ws.Range('E20').Value = '  +1.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.25'
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.48'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('E24').Value = '  +4.32%  '
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.28'
$ws.Range('E27').Value = '  +2.05%  '
$ws.Range('D28').Value = '0.0₃0770'
$ws.Range('E28').Value = '  +4.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  +3.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('E31').Value = '  +3.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.61'
$ws.Range('E32').Value = '  +4.72%  '
$ws.Range('E33').Value = '  +1.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.06'
$ws.Range('E34').Value = '  +4.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.914'
$ws.Range('E35').Value = '  +10.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.906'
$ws.Range('E36').Value = '  +11.26%  '
$ws.Range('E37').Value = '  +5.35%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.47'
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.50'
$ws.Range('E39').Value = '  +5.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '299.60'
$ws.Range('E40').Value = '  +6.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.65'
$ws.Range('E41').Value = '  +1.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.997'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.603'
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('E44').Value = '  +3.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0546'
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.31'
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.62'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('E48').Value = '  +5.23%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.09'
$ws.Range('E49').Value = '  +9.20%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.952.84'
$ws.Range('E50').Value = '  +0.58%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.58'
$ws.Range('E51').Value = '  +2.70%  '
